# repull data, push all data, mean calculation
# Update the F column (dSF) values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -8
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -4
